$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capacitors: package shrunk 0603 -> 0402, re-sourced to new manufacturers/parts ---
# C1, C4, C5 - 100nF
$ws.Range("D2").Value = "'0402_CAP"
$ws.Range("E2").Value = "'Wurth Electronics"
$ws.Range("F2").Value = "'885012105016"

# C2, C3 - 18pF
$ws.Range("D3").Value = "'0402_CAP"
$ws.Range("E3").Value = "'Walsin Technologies"
$ws.Range("F3").Value = "'0402N180F500CT"

# C6, C7 - 33pF
$ws.Range("D4").Value = "'0402_CAP"
$ws.Range("E4").Value = "'Wurth Electronics"
$ws.Range("F4").Value = "'885012005058"

# C8, C10 - 1uF/MLCC
$ws.Range("D5").Value = "'0402_CAP"
$ws.Range("E5").Value = "'Taiyo Yuden"
$ws.Range("F5").Value = "'JMK105BJ105KP-F"

# --- Resistors: package shrunk 0603 -> 0402, re-sourced to Vishay ---
# R1, R3, R4, R9, R14 - 10K
$ws.Range("D12").Value = "'0402_res"
$ws.Range("E12").Value = "'Vishay"
$ws.Range("F12").Value = "'CRCW040210K0FKEDC"

# R2 - 680R
$ws.Range("D13").Value = "'0402_res"
$ws.Range("E13").Value = "'Vishay"
$ws.Range("F13").Value = "'CRCW0402680RFKEDC"

# R5, R7 - 470R
$ws.Range("D14").Value = "'0402_RES"
$ws.Range("E14").Value = "'Vishay Semiconductors"
$ws.Range("F14").Value = "'CRCW0402470RFKEDC"

# R6, R8 - 33R
$ws.Range("D15").Value = "'0402_res"
$ws.Range("E15").Value = "'Vishay Semiconductors"
$ws.Range("F15").Value = "'CRCW040233R0FKEDC"

# R10 - 1.5K
$ws.Range("D16").Value = "'0402_res"
$ws.Range("E16").Value = "'Vishay"
$ws.Range("F16").Value = "'CRCW04021K50FKEDC"

# R12 - 0R
$ws.Range("D17").Value = "'0402_res"
$ws.Range("E17").Value = "'Vishay"
$ws.Range("F17").Value = "'CRCW04020000Z0EDC"

# R15 - 18K
$ws.Range("D18").Value = "'0402_RES"
$ws.Range("E18").Value = "'Vishay"
$ws.Range("F18").Value = "'CRCW040218K0FKED"

# --- Footer: time this BOM/basket was generated ---
$ws.Range("E24").Value = "'11:49"

# --- Row heights grow to fit the new, longer (wrapped) part numbers ---
$ws.Rows(3).RowHeight = 21
$ws.Rows(14).RowHeight = 21
$ws.Rows(15).RowHeight = 21
